$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 12-17 (entered first, top-to-bottom)
$ws.Range("C12").Value = "product_cost"
$ws.Range("C13").Value = "price_before_tax"
$ws.Range("C14").Value = "tax"
$ws.Range("C15").Value = "quantity"
$ws.Range("C16").Value = "discount"
$ws.Range("C17").Value = "price"

# Column E, rows 12-17
$ws.Range("E12").Value = "product_margin"
$ws.Range("E13").Value = "base_price"
$ws.Range("E17").Value = "base_price"

# Row 24
$ws.Range("C24").Value = "quantity"
$ws.Range("D24").Value = "price_before_tax"
$ws.Range("E24").Value = "tax"
$ws.Range("F24").Value = "price"

# Row 25
$ws.Range("D25").Value = "base_price"
$ws.Range("E25").Value = "per item"

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 183.75
$ws.Range("E26").Value = 0.05
$ws.Range("F26").Formula = "=D26+(D26*E26)"

# Row 27
$ws.Range("E27").Value = "c"

$ws.Range("F29").Select()
